# Update countries & provincias Spain
# - Belice moves ahead of Nueva Caledonia in the country list (row 192/193
#   swap names, each keeping its own stats).
# - Refreshed COVID case counters for several countries (rows 4, 10, 16,
#   33, 61, 144) and for the now-reordered Belice/Nueva Caledonia rows
#   (192, 193).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap country names for rows 192/193 -----------------------------
# Row 192 was "Nueva Caledonia"; it becomes "Belice".
# Row 193 was "Belice"; it becomes "Nueva Caledonia".
$ws.Range("A192").Value = "Belice"
$ws.Range("A193").Value = "Nueva Caledonia"

# --- Row 4 (Estados Unidos) -------------------------------------------
$ws.Range("B4").Value = 1376849
$ws.Range("C4").Value = 9211
$ws.Range("D4").Value = 258419
$ws.Range("E4").Value = 1037248
$ws.Range("G4").Value = 395
$ws.Range("H4").Value = 81182

# --- Row 10 (Alemania) -------------------------------------------------
$ws.Range("B10").Value = 172295
$ws.Range("C10").Value = 416
$ws.Range("E10").Value = 19084
$ws.Range("G10").Value = 42
$ws.Range("H10").Value = 7611

# --- Row 16 (Canada) -----------------------------------------------------
$ws.Range("B16").Value = 69907
$ws.Range("C16").Value = 1059
$ws.Range("D16").Value = 32650
$ws.Range("E16").Value = 32265
$ws.Range("G16").Value = 122
$ws.Range("H16").Value = 4992

# --- Row 33 (Israel) -----------------------------------------------------
$ws.Range("B33").Value = 16506
$ws.Range("C33").Value = 29
$ws.Range("D33").Value = 11843
$ws.Range("E33").Value = 4405
$ws.Range("G33").Value = 6
$ws.Range("H33").Value = 258

# --- Row 61 (Moldavia) ---------------------------------------------------
$ws.Range("B61").Value = 4995
$ws.Range("C61").Value = 68
$ws.Range("E61").Value = 2840
$ws.Range("G61").Value = 6
$ws.Range("H61").Value = 175

# --- Row 144 (Madagascar) -------------------------------------------------
$ws.Range("C144").Value = 7

# --- Row 192 (now Belice, after swap) -------------------------------------
$ws.Range("D192").Value = 16
$ws.Range("H192").Value = 2

# --- Row 193 (now Nueva Caledonia, after swap) ----------------------------
$ws.Range("D193").Value = 18
$ws.Range("H193").Value = 0
